# Applies the 21-12-2023 02:45 script update to the LaLiga2 2023-2024 sheet:
#  - several adjacent match rows were re-sorted by kickoff time, so their
#    F:V (home..url) payloads swap/rotate while A:E (index/meta/date) stay put
#  - three newly scraped matches are appended at the bottom (rows 226-228)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sorted row pairs: swap the F:V match payload between the two rows ---
$pair49_50_top = $ws.Range("F49:V49").Value()
$pair49_50_bot = $ws.Range("F50:V50").Value()
$ws.Range("F49:V49").Value = $pair49_50_bot
$ws.Range("F50:V50").Value = $pair49_50_top

$pair59_60_top = $ws.Range("F59:V59").Value()
$pair59_60_bot = $ws.Range("F60:V60").Value()
$ws.Range("F59:V59").Value = $pair59_60_bot
$ws.Range("F60:V60").Value = $pair59_60_top

$pair97_98_top = $ws.Range("F97:V97").Value()
$pair97_98_bot = $ws.Range("F98:V98").Value()
$ws.Range("F97:V97").Value = $pair97_98_bot
$ws.Range("F98:V98").Value = $pair97_98_top

$pair99_100_top = $ws.Range("F99:V99").Value()
$pair99_100_bot = $ws.Range("F100:V100").Value()
$ws.Range("F99:V99").Value = $pair99_100_bot
$ws.Range("F100:V100").Value = $pair99_100_top

$pair114_115_top = $ws.Range("F114:V114").Value()
$pair114_115_bot = $ws.Range("F115:V115").Value()
$ws.Range("F114:V114").Value = $pair114_115_bot
$ws.Range("F115:V115").Value = $pair114_115_top

$pair124_125_top = $ws.Range("F124:V124").Value()
$pair124_125_bot = $ws.Range("F125:V125").Value()
$ws.Range("F124:V124").Value = $pair124_125_bot
$ws.Range("F125:V125").Value = $pair124_125_top

$pair169_170_top = $ws.Range("F169:V169").Value()
$pair169_170_bot = $ws.Range("F170:V170").Value()
$ws.Range("F169:V169").Value = $pair169_170_bot
$ws.Range("F170:V170").Value = $pair169_170_top

$pair179_180_top = $ws.Range("F179:V179").Value()
$pair179_180_bot = $ws.Range("F180:V180").Value()
$ws.Range("F179:V179").Value = $pair179_180_bot
$ws.Range("F180:V180").Value = $pair179_180_top

$pair184_185_top = $ws.Range("F184:V184").Value()
$pair184_185_bot = $ws.Range("F185:V185").Value()
$ws.Range("F184:V184").Value = $pair184_185_bot
$ws.Range("F185:V185").Value = $pair184_185_top

$pair196_197_top = $ws.Range("F196:V196").Value()
$pair196_197_bot = $ws.Range("F197:V197").Value()
$ws.Range("F196:V196").Value = $pair196_197_bot
$ws.Range("F197:V197").Value = $pair196_197_top

# --- Re-sorted 3-row block 85-87: rotate the F:V match payload (87->85->86->87) ---
$row85 = $ws.Range("F85:V85").Value()
$row86 = $ws.Range("F86:V86").Value()
$row87 = $ws.Range("F87:V87").Value()
$ws.Range("F85:V85").Value = $row87
$ws.Range("F86:V86").Value = $row85
$ws.Range("F87:V87").Value = $row86

# --- Newly scraped matches appended at the bottom ---
$ws.Cells.Item(226,1).Value = 225
$ws.Cells.Item(226,2).Value = "spain"
$ws.Cells.Item(226,3).Value = "laliga2"
$ws.Cells.Item(226,4).Value = "2023-2024"
$ws.Cells.Item(226,5).Value = 45280.79166666666
$ws.Cells.Item(226,6).Value = "Zaragoza"
$ws.Cells.Item(226,7).Value = 2
$ws.Cells.Item(226,8).Value = "Levante"
$ws.Cells.Item(226,9).Value = 2
$ws.Cells.Item(226,10).Value = 2.24
$ws.Cells.Item(226,11).Value = "17/12/2023 18:43"
$ws.Cells.Item(226,12).Value = 2.25
$ws.Cells.Item(226,13).Value = "20/12/2023 18:58"
$ws.Cells.Item(226,14).Value = 3.08
$ws.Cells.Item(226,15).Value = "17/12/2023 18:43"
$ws.Cells.Item(226,16).Value = 3.07
$ws.Cells.Item(226,17).Value = "20/12/2023 18:58"
$ws.Cells.Item(226,18).Value = 3.73
$ws.Cells.Item(226,19).Value = "17/12/2023 18:43"
$ws.Cells.Item(226,20).Value = 3.88
$ws.Cells.Item(226,21).Value = "20/12/2023 18:58"
$ws.Cells.Item(226,22).Value = "https://www.betexplorer.com/football/spain/laliga2/zaragoza-levante/OxlEugH5/"

$ws.Cells.Item(227,1).Value = 226
$ws.Cells.Item(227,2).Value = "spain"
$ws.Cells.Item(227,3).Value = "laliga2"
$ws.Cells.Item(227,4).Value = "2023-2024"
$ws.Cells.Item(227,5).Value = 45280.89583333334
$ws.Cells.Item(227,6).Value = "Eibar"
$ws.Cells.Item(227,7).Value = 1
$ws.Cells.Item(227,8).Value = "Gijon"
$ws.Cells.Item(227,9).Value = 1
$ws.Cells.Item(227,10).Value = 1.89
$ws.Cells.Item(227,11).Value = "17/12/2023 14:12"
$ws.Cells.Item(227,12).Value = 1.98
$ws.Cells.Item(227,13).Value = "20/12/2023 21:29"
$ws.Cells.Item(227,14).Value = 3.43
$ws.Cells.Item(227,15).Value = "17/12/2023 14:12"
$ws.Cells.Item(227,16).Value = 3.39
$ws.Cells.Item(227,17).Value = "20/12/2023 21:29"
$ws.Cells.Item(227,18).Value = 4.53
$ws.Cells.Item(227,19).Value = "17/12/2023 14:12"
$ws.Cells.Item(227,20).Value = 4.37
$ws.Cells.Item(227,21).Value = "20/12/2023 21:29"
$ws.Cells.Item(227,22).Value = "https://www.betexplorer.com/football/spain/laliga2/eibar-gijon/v9SrQDHI/"

$ws.Cells.Item(228,1).Value = 227
$ws.Cells.Item(228,2).Value = "spain"
$ws.Cells.Item(228,3).Value = "laliga2"
$ws.Cells.Item(228,4).Value = "2023-2024"
$ws.Cells.Item(228,5).Value = 45280.89583333334
$ws.Cells.Item(228,6).Value = "Leganes"
$ws.Cells.Item(228,7).Value = 1
$ws.Cells.Item(228,8).Value = "Tenerife"
$ws.Cells.Item(228,9).Value = 1
$ws.Cells.Item(228,10).Value = 2.27
$ws.Cells.Item(228,11).Value = "17/12/2023 21:12"
$ws.Cells.Item(228,12).Value = 2.44
$ws.Cells.Item(228,13).Value = "20/12/2023 21:16"
$ws.Cells.Item(228,14).Value = 2.92
$ws.Cells.Item(228,15).Value = "17/12/2023 21:12"
$ws.Cells.Item(228,16).Value = 2.69
$ws.Cells.Item(228,17).Value = "20/12/2023 21:16"
$ws.Cells.Item(228,18).Value = 3.9
$ws.Cells.Item(228,19).Value = "17/12/2023 21:12"
$ws.Cells.Item(228,20).Value = 4.06
$ws.Cells.Item(228,21).Value = "20/12/2023 21:16"
$ws.Cells.Item(228,22).Value = "https://www.betexplorer.com/football/spain/laliga2/leganes-tenerife/K45VmH9P/"

# Match the existing "Indice" (col A, bold/boxed) and "data_partida" (col E, date) styling
# used by every other row in the sheet, instead of leaving the new rows unstyled.
$ws.Range("A2").Copy()
$ws.Range("A226").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E226").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A227").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E227").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A228").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E228").PasteSpecial(-4122)

$excel.CutCopyMode = 0
